$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22 (shifting FLAMOGEST and everything below it down by one row)
$ws.Rows.Item(22).Insert()

# Populate the new row 22 with the FAROVIGA item data
$ws.Cells.Item(22, 1).Value = 16
$ws.Cells.Item(22, 3).Value = "FAROVIGA 100MG 12 F.C.TAB."
$ws.Cells.Item(22, 8).Value = "6:0"
$ws.Cells.Item(22, 12).Value = "1"
$ws.Cells.Item(22, 14).Value = "108.00"
$ws.Cells.Item(22, 16).Value = "17.2800"
$ws.Cells.Item(22, 17).Value = "0:2"

# Update the subtotal (now row 48) and the generated timestamp
$ws.Cells.Item(48, 16).Value = 2527.5
$ws.Cells.Item(49, 1).Value = "Sunday, 5 October, 2025 1:08 PM"
